$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended after the 2025-07-13 entry, for the 2025-07-16 run.
# Column A holds the date as literal text (matching the existing rows' format),
# so force a text number format before assigning the string value, then reset
# the cell style back to Normal so no stray formatting is left behind.
$ws.Cells.Item(32, 1).NumberFormat = "@"
$ws.Cells.Item(32, 1).Value = "07/16/2025"
$ws.Cells.Item(32, 1).Style = "Normal"

$ws.Cells.Item(32, 2).Value = 0.0004174600000000014
$ws.Cells.Item(32, 3).Value = 119771.9541992043
$ws.Cells.Item(32, 4).Value = 50
